$d = $word.ActiveDocument

$replacements = @(
    @("70÷2=", "91÷3="),
    @("75÷9=", "68÷6="),
    @("20÷9=", "54÷2="),
    @("10÷6=", "99÷9="),
    @("41÷4=", "99÷5="),
    @("37÷4=", "58÷7="),
    @("79÷6=", "64÷7="),
    @("65÷8=", "94÷8="),
    @("46÷7=", "75÷2="),
    @("49÷7=", "13÷5="),
    @("52÷9=", "13÷6="),
    @("28÷7=", "16÷4="),
    @("86÷8=", "21÷5="),
    @("51÷7=", "44÷2="),
    @("93÷7=", "31÷6="),
    @("50÷5=", "29÷4="),
    @("23÷5=", "22÷6="),
    @("48÷8=", "20÷3="),
    @("34÷6=", "27÷5="),
    @("75÷6=", "17÷4="),
    @("81÷8=", "38÷8="),
    @("95÷2=", "20÷2="),
    @("39÷7=", "31÷4="),
    @("24÷6=", "11÷8="),
    @("95÷3=", "17÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
